$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: new changelog entry for 0.5.4
$ws.Range("A19").Value2 = 46060
$ws.Range("A19").NumberFormat = "yyyy-mm-dd"
$ws.Range("B19").Value2 = "0.5.4"
$ws.Range("C19").Value2 = "Enhancement"
$ws.Range("D19").Value2 = "ListPage component with built-in search and FilterBar toolbar. Standardised filtering across all pages using FilterBar (replaces standalone Select dropdowns). Migrated training courses page. Documented UI layout conventions (work/ui-layout/)."
$ws.Range("E19").Value2 = "Claude"

# Row 20: new changelog entry for 0.5.5
$ws.Range("A20").Value2 = 46060
$ws.Range("A20").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B20").Value2 = "0.5.5"
$ws.Range("C20").Value2 = "Enhancement"
$ws.Range("D20").Value2 = "DataTable component for tabular records (flat or expandable rows). Merged notifications into audit log page as tab. Migrated all list pages to ListPage+ListRow+FilterBar: employees, users, my-training, leave, onboarding items, employee statuses. Updated employee-leave-tab to use ListRow."
$ws.Range("E20").Value2 = "Claude"
